$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '30.514.30'
Set-TextValue $ws.Range('E2') '  +1.37%  '
Set-TextValue $ws.Range('D3') '1.879.04'
Set-TextValue $ws.Range('E3') '  +1.17%  '
Set-TextValue $ws.Range('E4') '  +0.02%  '
Set-TextValue $ws.Range('D5') '247.32'
Set-TextValue $ws.Range('E5') '  +6.00%  '
Set-TextValue $ws.Range('E6') '  -0.02%  '
Set-TextValue $ws.Range('D7') '0.4745'
Set-TextValue $ws.Range('E7') '  +1.64%  '
Set-TextValue $ws.Range('D8') '0.2916'
Set-TextValue $ws.Range('E8') '  +3.05%  '
Set-TextValue $ws.Range('D9') '0.06519'
Set-TextValue $ws.Range('E9') '  +1.17%  '
Set-TextValue $ws.Range('D10') '22.00'
Set-TextValue $ws.Range('E10') '  +6.36%  '
Set-TextValue $ws.Range('D11') '97.94'
Set-TextValue $ws.Range('E11') '  +4.97%  '
Set-TextValue $ws.Range('D12') '0.07722'
Set-TextValue $ws.Range('E12') '  +0.76%  '
Set-TextValue $ws.Range('D13') '0.7391'
Set-TextValue $ws.Range('E13') '  +9.40%  '
Set-TextValue $ws.Range('D14') '1.880.56'
Set-TextValue $ws.Range('E14') '  +1.69%  '
Set-TextValue $ws.Range('D15') '5.132'
Set-TextValue $ws.Range('E15') '  +2.07%  '
Set-TextValue $ws.Range('D16') '273.82'
Set-TextValue $ws.Range('E16') '  +2.63%  '
Set-TextValue $ws.Range('D17') '30.508.34'
Set-TextValue $ws.Range('E17') '  +1.44%  '
Set-TextValue $ws.Range('D18') '13.44'
Set-TextValue $ws.Range('E18') '  +1.31%  '
Set-TextValue $ws.Range('D19') '0.000007559'
Set-TextValue $ws.Range('E19') '  +1.10%  '
Set-TextValue $ws.Range('D20') '1.001'
Set-TextValue $ws.Range('E20') '  +0.03%  '
Set-TextValue $ws.Range('D21') '2.130.58'
Set-TextValue $ws.Range('E21') '  +1.91%  '
Set-TextValue $ws.Range('E22') '  +0.07%  '
Set-TextValue $ws.Range('D23') '5.245'
Set-TextValue $ws.Range('E23') '  +2.31%  '
Set-TextValue $ws.Range('E24') '  +1.98%  '
Set-TextValue $ws.Range('D25') '9.280'
Set-TextValue $ws.Range('E25') '  +0.39%  '
Set-TextValue $ws.Range('D26') '163.59'
Set-TextValue $ws.Range('E26') '  -1.03%  '
Set-TextValue $ws.Range('D27') '18.88'
Set-TextValue $ws.Range('E27') '  +1.71%  '
Set-TextValue $ws.Range('D28') '1.941'
Set-TextValue $ws.Range('E28') '  +3.77%  '
Set-TextValue $ws.Range('E29') '  +3.21%  '
Set-TextValue $ws.Range('D30') '1.368'
Set-TextValue $ws.Range('E30') '  -0.31%  '
Set-TextValue $ws.Range('D31') '1.514'
Set-TextValue $ws.Range('E31') '  +4.10%  '
Set-TextValue $ws.Range('D32') '4.322'
Set-TextValue $ws.Range('E32') '  +3.52%  '
Set-TextValue $ws.Range('D33') '4.100'
Set-TextValue $ws.Range('E33') '  +3.60%  '
Set-TextValue $ws.Range('D34') '0.04823'
Set-TextValue $ws.Range('E34') '  +4.22%  '
Set-TextValue $ws.Range('D35') '1.127'
Set-TextValue $ws.Range('E35') '  +2.07%  '
Set-TextValue $ws.Range('D36') '0.7009'
Set-TextValue $ws.Range('E36') '  +3.18%  '
Set-TextValue $ws.Range('B37') 'Frax'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range('D37') '1.000'
Set-TextValue $ws.Range('E37') '  -0.01%  '
Set-TextValue $ws.Range('B38') 'HuobiToken'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D38') '2.718'
Set-TextValue $ws.Range('E38') '  +0.22%  '
Set-TextValue $ws.Range('B39') 'VeChain'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D39') '0.01864'
Set-TextValue $ws.Range('E39') '  +2.68%  '
Set-TextValue $ws.Range('B40') 'MXToken'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D40') '2.748'
Set-TextValue $ws.Range('E40') '  +1.53%  '
Set-TextValue $ws.Range('B41') 'FraxShare'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D41') '6.314'
Set-TextValue $ws.Range('E41') '  +0.75%  '
Set-TextValue $ws.Range('B42') 'RenderToken'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D42') '1.980'
Set-TextValue $ws.Range('E42') '  +6.22%  '
Set-TextValue $ws.Range('B43') 'Aave'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D43') '71.08'
Set-TextValue $ws.Range('E43') '  +1.83%  '
Set-TextValue $ws.Range('B44') 'TheSandbox'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D44') '0.4211'
Set-TextValue $ws.Range('E44') '  +4.74%  '
Set-TextValue $ws.Range('E45') '  +0.00%  '
Set-TextValue $ws.Range('B46') 'TrustWalletToken'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D46') '0.8386'
Set-TextValue $ws.Range('E46') '  +1.76%  '
Set-TextValue $ws.Range('B47') 'Quant'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D47') '102.61'
Set-TextValue $ws.Range('E47') '  +0.79%  '
Set-TextValue $ws.Range('B48') 'EnergySwap'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D48') '9.326'
Set-TextValue $ws.Range('E48') '  +2.58%  '
Set-TextValue $ws.Range('B49') 'Aptos'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D49') '7.079'
Set-TextValue $ws.Range('E49') '  +2.94%  '
Set-TextValue $ws.Range('B50') 'Elrond'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws.Range('D50') '35.61'
Set-TextValue $ws.Range('E50') '  +4.92%  '
Set-TextValue $ws.Range('B51') 'Maker'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D51') '915.87'
Set-TextValue $ws.Range('E51') '  -0.38%  '
